$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ideal-format")

# Replace static N values (row 3) with formulas summing the S,E,I_asym,I_sym,I_sev,R,D rows (16-22)
$ws.Range("B3").Formula = "=B16+B17+B18+B19+B20+B21+B22"
$ws.Range("C3").Formula = "=C16+C17+C18+C19+C20+C21+C22"

# Move the active selection to C3 (matches the author's final cursor position)
$ws.Range("C3").Select()
